$d = $word.ActiveDocument

# --- Table 1: "Nom de balise" header table (caseId / mobilizedResource) ---
$tbl1 = $d.Tables.Item(1)
$rowMobilized = $tbl1.Rows.Item(3)

# 1) "Ressource engagée" -> "Ressource engagée / à engager"
$rowMobilized.Cells.Item(2).Range.Text = "Ressource engagée / à engager"

# 2) Description cell: append extra sentence + line break + new sentence
$rowMobilized.Cells.Item(5).Range.Text = "Objet permettant de communquer la liste des ressource et vecteurs mobilisés en 15-15 et 15-SMUR pour le message RS-RI" + [char]11 + "Objet permettant de communiquer la liste des ressources à engager en 15-SMUR pour le message RS-ER"

# --- Table 2: "resource" object table ---
$tbl2 = $d.Tables.Item(2)

# 3) Remove the whole "resourceType" row
for ($i = 1; $i -le $tbl2.Rows.Count; $i++) {
    $row = $tbl2.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text.TrimEnd([char]7, [char]13) -eq "resourceType") {
        $row.Delete()
        break
    }
}

# 4) "vehiculeType" row cardinality: 0..1 -> 1..1
for ($i = 1; $i -le $tbl2.Rows.Count; $i++) {
    $row = $tbl2.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text.TrimEnd([char]7, [char]13) -eq "vehiculeType") {
        $row.Cells.Item(4).Range.Text = "1..1"
        break
    }
}

# 5) "Etats vecteur" -> "Etat vecteur"
for ($i = 1; $i -le $tbl2.Rows.Count; $i++) {
    $row = $tbl2.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text.TrimEnd([char]7, [char]13) -eq "state") {
        $row.Cells.Item(2).Range.Text = "Etat vecteur"
        $row.Cells.Item(5).Range.Text = "Objet qui permet de décrire l'état d'un vecteur mobilisé - sous forme de liste, il permet de décrire l'historique des états connus d'un même vecteur. "
        break
    }
}
